# "Input_Value" sheet holds the test-data credentials used by the Selenium
# automation (URL / UserName / Password) in M2:O2. The uploaded workbook had
# these sample credentials scrubbed out before being committed, which is why
# the corresponding shared strings (the URL, "IBM_IMPLEMENTATION_USER" and
# "Oracle1234") disappear from the package once the cells no longer
# reference them.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

$ws.Range("M2:O2").ClearContents() | Out-Null

# Leave the sheet showing that range selected, approximating the author's
# on-screen selection when the file was saved.
$ws.Activate()
$ws.Range("M2:O2").Select() | Out-Null
